# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Re-sorts the worker/period detail table (rows 16-21) on sheet "Hoja1"
# by period ascending (column E), inserting LUIS SEGUNDO PASTRANA
# HERNANDEZ's single 1808 record among JUAN DAVID MADRID OCHOA's
# 1806-1810 records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order/content for the detail table B16:G21
# (Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico)
$data = @(
    @("CC", "1068391397", "JUAN DAVID MADRID OCHOA",        "1806", 27083, 781242),
    @("CC", "1068391397", "JUAN DAVID MADRID OCHOA",        "1807", 31249, 781242),
    @("CC", "70526895",   "LUIS SEGUNDO PASTRANA HERNANDEZ", "1808", 42000, 1050000),
    @("CC", "1068391397", "JUAN DAVID MADRID OCHOA",        "1808", 31249, 781242),
    @("CC", "1068391397", "JUAN DAVID MADRID OCHOA",        "1809", 31249, 781242),
    @("CC", "1068391397", "JUAN DAVID MADRID OCHOA",        "1810", 31249, 781242)
)

$startRow = 16
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
